$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 209.125
$ws.Range("I53").Value = 284.44446
$ws.Range("K53").Value = 284.44446
$ws.Range("M53").Value = 352.55554
# Row 113
$ws.Range("H113").Value = 4546
$ws.Range("I113").Value = 3365
$ws.Range("K113").Value = 3365
$ws.Range("M113").Value = -111
# Row 138
$ws.Range("H138").Value = 5128.7617
$ws.Range("I138").Value = 4899
$ws.Range("K138").Value = 14697
$ws.Range("M138").Value = -9557
# Row 141
$ws.Range("H141").Value = 62382
$ws.Range("I141").Value = 66623.57000000001
$ws.Range("K141").Value = 199870.71
$ws.Range("M141").Value = -194690.71

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 333607.16
$ws.Range("I4").Value = 333607.16
$ws.Range("K4").Value = 333607.16
$ws.Range("M4").Value = -333491.16
# Row 37
$ws.Range("H37").Value = 25015000
$ws.Range("J37").Value = 30000
$ws.Range("L37").Value = 30000
$ws.Range("N37").Value = -30546
# Row 45
$ws.Range("H45").Value = 33076
$ws.Range("I45").Value = 35791.832
$ws.Range("J45").Value = 29002.25
$ws.Range("K45").Value = 35791.832
$ws.Range("L45").Value = 29002.25
$ws.Range("M45").Value = -35414.832
$ws.Range("N45").Value = -29756.25
# Row 110
$ws.Range("H110").Value = 2120
$ws.Range("I110").Value = 2120
$ws.Range("K110").Value = 2120
$ws.Range("M110").Value = -75
# Row 122
$ws.Range("H122").Value = 13515.883
$ws.Range("I122").Value = 17104.309
$ws.Range("K122").Value = 51312.927
$ws.Range("M122").Value = -48862.927
# Row 132
$ws.Range("H132").Value = 2652.4285
$ws.Range("I132").Value = 2515.375
$ws.Range("K132").Value = 7546.125
$ws.Range("M132").Value = -5016.125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 87
$ws.Range("H87").Value = 30354
$ws.Range("J87").Value = 30354
$ws.Range("L87").Value = 30354
$ws.Range("N87").Value = -32850
# Row 90
$ws.Range("H90").Value = 30354
$ws.Range("J90").Value = 30354
$ws.Range("L90").Value = 91062
$ws.Range("N90").Value = -103542
# Row 107
$ws.Range("H107").Value = 1868.125
$ws.Range("I107").Value = 1157.5
$ws.Range("K107").Value = 1157.5
$ws.Range("M107").Value = 762.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 10001873
$ws.Range("I4").Value = 2497.3333
$ws.Range("K4").Value = 2497.3333
$ws.Range("M4").Value = -2385.3333
# Row 37
$ws.Range("H37").Value = 15000
$ws.Range("J37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15214
# Row 41
$ws.Range("H41").Value = 8136
$ws.Range("J41").Value = 8368
$ws.Range("L41").Value = 8368
$ws.Range("N41").Value = -9224
# Row 59
$ws.Range("H59").Value = 17495
$ws.Range("J59").Value = 17495
$ws.Range("L59").Value = 17495
$ws.Range("N59").Value = -19785
# Row 68
$ws.Range("H68").Value = 24997.777
$ws.Range("J68").Value = 24997.777
$ws.Range("L68").Value = 24997.777
$ws.Range("N68").Value = -26495.777
# Row 71
$ws.Range("H71").Value = 24997.777
$ws.Range("J71").Value = 24997.777
$ws.Range("L71").Value = 74993.33099999999
$ws.Range("N71").Value = -82481.33099999999
# Row 74
$ws.Range("H74").Value = 35500
$ws.Range("J74").Value = 35500
$ws.Range("L74").Value = 35500
$ws.Range("N74").Value = -37248
# Row 77
$ws.Range("H77").Value = 35500
$ws.Range("J77").Value = 35500
$ws.Range("L77").Value = 106500
$ws.Range("N77").Value = -115236
# Row 99
$ws.Range("H99").Value = 9800
$ws.Range("I99").Value = 9000
$ws.Range("K99").Value = 9000
$ws.Range("M99").Value = -7502
# Row 122
$ws.Range("H122").Value = 1887.2916
$ws.Range("J122").Value = 4020
$ws.Range("L122").Value = 12060
$ws.Range("N122").Value = -16960
# Row 124
$ws.Range("H124").Value = 72250.25
$ws.Range("J124").Value = 72250.25
$ws.Range("L124").Value = 72250.25
$ws.Range("N124").Value = -77160.25
# Row 126
$ws.Range("H126").Value = 9800
$ws.Range("I126").Value = 9000
$ws.Range("K126").Value = 27000
$ws.Range("M126").Value = -24530

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 87.5
$ws.Range("I2").Value = 79
$ws.Range("J2").Value = 130
$ws.Range("K2").Value = 474
$ws.Range("L2").Value = 780
$ws.Range("M2").Value = -361
$ws.Range("N2").Value = -1006
# Row 107
$ws.Range("H107").Value = 1371.4286
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 1361.5385
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 4084.6155
$ws.Range("M107").Value = -2580
$ws.Range("N107").Value = -7924.6155
# Row 124
$ws.Range("H124").Value = 18900
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 18900
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 56700
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -66520
# Row 126
$ws.Range("H126").Value = 19444
$ws.Range("I126").Value = 8988
$ws.Range("J126").Value = 29900
$ws.Range("K126").Value = 26964
$ws.Range("L126").Value = 89700
$ws.Range("M126").Value = -22024
$ws.Range("N126").Value = -99580
# Row 128
$ws.Range("H128").Value = 309318.78
$ws.Range("I128").Value = 309318.78
$ws.Range("K128").Value = 927956.3400000001
$ws.Range("M128").Value = -922976.3400000001
# Row 133
$ws.Range("H133").Value = 6815.7856
$ws.Range("I133").Value = 6722.231
$ws.Range("J133").Value = 8032
$ws.Range("K133").Value = 20166.693
$ws.Range("L133").Value = 24096
$ws.Range("M133").Value = -15106.693
$ws.Range("N133").Value = -34216
# Row 139
$ws.Range("H139").Value = 4895.2
$ws.Range("I139").Value = 6493.3335
$ws.Range("J139").Value = 2498
$ws.Range("K139").Value = 19480.0005
$ws.Range("L139").Value = 7494
$ws.Range("M139").Value = -14340.0005
$ws.Range("N139").Value = -17774

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 123
$ws.Range("H123").Value = 46035
$ws.Range("J123").Value = 46035
$ws.Range("L123").Value = 46035
$ws.Range("N123").Value = -50935
# Row 132
$ws.Range("H132").Value = 2165.111
$ws.Range("I132").Value = 2369.4285
$ws.Range("K132").Value = 7108.2855
$ws.Range("M132").Value = -4578.2855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1632.3334
$ws.Range("I22").Value = 1468.8
$ws.Range("J22").Value = 2450
$ws.Range("K22").Value = 1468.8
$ws.Range("L22").Value = 2450
$ws.Range("M22").Value = -1173.8
$ws.Range("N22").Value = -3040
# Row 27
$ws.Range("H27").Value = 1632.3334
$ws.Range("I27").Value = 1468.8
$ws.Range("J27").Value = 2450
$ws.Range("K27").Value = 1468.8
$ws.Range("L27").Value = 2450
$ws.Range("M27").Value = -1361.8
$ws.Range("N27").Value = -2664
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 93
$ws.Range("H93").Value = 48960.285
$ws.Range("I93").Value = 1519.8
$ws.Range("K93").Value = 1519.8
$ws.Range("M93").Value = -271.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
